# aggiornamento fino a 02/05
# Append 6 new daily rows (239-244) to Sheet1, continuing the existing
# date/nuovi-pos./somma-mobile-7gg./somma-mobile-per-100k series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$newRows = @(
    @(44313, 0, 13, 131.8191036300953),
    @(44314, 1, 13, 131.8191036300953),
    @(44315, 2, 14, 141.9590346785642),
    @(44316, 2, 8,  81.11944838775096),
    @(44317, 2, 8,  81.11944838775096),
    @(44318, 2, 9,  91.25937943621983)
)

$lastExistingRow = 238
$startRow = $lastExistingRow + 1

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]

    # Copy the date cell's formatting (border/font/alignment/date numfmt)
    # from the last existing row so the new date cell matches the
    # existing column-A style instead of picking up a plain default.
    $ws.Range("A$lastExistingRow").Copy($ws.Range("A$r"))

    $ws.Range("A$r").Value = $values[0]
    $ws.Range("B$r").Value = $values[1]
    $ws.Range("C$r").Value = $values[2]
    $ws.Range("D$r").Value = $values[3]
}
